$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 9 to make room for "MANAGEMENT NETMASK" - shifts
# MANAGEMENT_GATEWAY..DEFAULT_ROUTE and the VLAN table below down by one row.
$ws.Rows(9).Insert()

# --- Firewall Basic Info table (rows 1-14) -------------------------------
$ws.Range("A1").Value = "FIELD NAME"
$ws.Range("B1").Value = "USER INPUT"

$ws.Range("A2").Value = "FIREWALL MODEL"
$ws.Range("B2").Value = "200F"

$ws.Range("A3").Value = "SITE CODE"
$ws.Range("B3").Value = "USABD"

$ws.Range("A4").Value = "CLUSTER NUMBER"
$ws.Range("B4").Value = 1

$ws.Range("A5").Value = "SEGMENTATION TYPE"
$ws.Range("B5").Value = "MFG"

$ws.Range("A6").Value = "REGION"
$ws.Range("B6").Value = "AMER"

$ws.Range("A7").Value = "MANAGEMENT INTERFACE"
$ws.Range("B7").Value = "mgmt"

$ws.Range("A8").Value = "MANAGEMENT IPS"
$ws.Range("B8").Value = "10.1.1.15, 10.1.1.16"

$ws.Range("A9").Value = "MANAGEMENT NETMASK"
$ws.Range("B9").Value = 24

$ws.Range("A10").Value = "MANAGEMENT GATEWAY"
$ws.Range("B10").Value = "10.1.1.1"

$ws.Range("A11").Value = "HA MODE"
$ws.Range("B11").Value = "ACTIVE-PASSIVE"

$ws.Range("A12").Value = "HA INTERFACES"
$ws.Range("B12").Value = "PORT1,PORT2"

$ws.Range("A13").Value = "LACP INTERFACES"
$ws.Range("B13").Value = "X1,X2"

$ws.Range("A14").Value = "DEFAULT ROUTE"
$ws.Range("B14").Value = "10.1.2.1"

# Bold + centered style for the field-name column of the basic info table.
$ws.Range("A2:A14").Font.Bold = $true
$ws.Range("A2:A14").HorizontalAlignment = -4108
$ws.Range("A2:A14").VerticalAlignment = -4108

# Move the active selection to reflect the edited cell.
$ws.Range("A14").Select()
